# Elimna EC anteriores y se agregan nuevos, se modifica base de datos
#
# The "Periodo Mora" column (E) for the JOSE VICTOR HERRERA TORRES block
# (rows 17-23) is re-sequenced to descending order (2407 -> 2401), and the
# "Valor Mora" amounts (F) for the first/last rows of that block are swapped
# so the 153334 remainder now lands on the newest period (2407) instead of
# the oldest (2401).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("E17").Value = "2407"
$ws.Range("F17").Value = 153334

$ws.Range("E18").Value = "2406"
$ws.Range("F18").Value = 200000

$ws.Range("E19").Value = "2405"
$ws.Range("F19").Value = 200000

$ws.Range("E20").Value = "2404"
$ws.Range("F20").Value = 200000

$ws.Range("E21").Value = "2403"
$ws.Range("F21").Value = 200000

$ws.Range("E22").Value = "2402"
$ws.Range("F22").Value = 200000

$ws.Range("E23").Value = "2401"
$ws.Range("F23").Value = 200000
